$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-13) across columns B-E (2-5) hold floating point values
# that must be rewritten to disk as rounded integers.
$startRow = 2
$endRow = 13
$startCol = 2
$endCol = 5

for ($r = $startRow; $r -le $endRow; $r++) {
    for ($c = $startCol; $c -le $endCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null) {
            $d = [double]$val
            if ($d -ge 0) {
                $rounded = [Math]::Floor($d + 0.5)
            } else {
                $rounded = [Math]::Ceiling($d - 0.5)
            }
            $cell.Value = $rounded
        }
    }
}
